$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Recid"
$ws.Range("B1").Value = "projId"
$ws.Range("C1").Value = "taskname"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "suporte Empresa 1"
$ws.Range("C2").Value = "desenvolvimento"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "suporte Empresa 1"
$ws.Range("C3").Value = "suporte"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "suporte Empresa 1"
$ws.Range("C4").Value = "reunião"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "suporte Empresa 1"
$ws.Range("C5").Value = "acompanhamento"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "suporte Empresa 1"
$ws.Range("C6").Value = "documentação"

# Remove the underline font style previously applied to B1/C1 (now none).
$ws.Range("B1:C1").Font.Underline = $false

# Clear prior per-cell selection/view tweaks are not directly scriptable via COM;
# Excel will persist current view state on save.
